$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 5551890
$ws.Range("I33").Value = 8327721
$ws.Range("K33").Value = 8327721
$ws.Range("M33").Value = -8327492
$ws.Range("H55").Value = 687.5
$ws.Range("I55").Value = 200
$ws.Range("J55").Value = 757.1429000000001
$ws.Range("K55").Value = 200
$ws.Range("L55").Value = 757.1429000000001
$ws.Range("M55").Value = 14
$ws.Range("N55").Value = -1185.1429
$ws.Range("H129").Value = 948.1579
$ws.Range("I129").Value = 404.66666
$ws.Range("J129").Value = 1050.0625
$ws.Range("K129").Value = 1213.99998
$ws.Range("L129").Value = 3150.1875
$ws.Range("M129").Value = 3786.00002
$ws.Range("N129").Value = -13150.1875

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 36000
$ws.Range("J24").Value = 36000
$ws.Range("L24").Value = 36000
$ws.Range("N24").Value = -36748
$ws.Range("H74").Value = 8927.444
$ws.Range("I74").Value = 3298.2727
$ws.Range("J74").Value = 17773.285
$ws.Range("K74").Value = 3298.2727
$ws.Range("L74").Value = 17773.285
$ws.Range("M74").Value = -2424.2727
$ws.Range("N74").Value = -19521.285
$ws.Range("H77").Value = 8927.444
$ws.Range("I77").Value = 3298.2727
$ws.Range("J77").Value = 17773.285
$ws.Range("K77").Value = 16491.3635
$ws.Range("L77").Value = 88866.425
$ws.Range("M77").Value = -12123.3635
$ws.Range("N77").Value = -97602.425
$ws.Range("H100").Value = 36000
$ws.Range("J100").Value = 36000
$ws.Range("L100").Value = 36000
$ws.Range("N100").Value = -38164
$ws.Range("H102").Value = 5294391
$ws.Range("I102").Value = 6176056
$ws.Range("K102").Value = 6176056
$ws.Range("M102").Value = -6174434
$ws.Range("H132").Value = 3226.5
$ws.Range("I132").Value = 2904.5454
$ws.Range("J132").Value = 4997.25
$ws.Range("K132").Value = 8713.636200000001
$ws.Range("L132").Value = 14991.75
$ws.Range("M132").Value = -6183.636200000001
$ws.Range("N132").Value = -20051.75
$ws.Range("H134").Value = 53331.668
$ws.Range("J134").Value = 53331.668
$ws.Range("L134").Value = 53331.668
$ws.Range("N134").Value = -63471.668

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1262.1111
$ws.Range("I99").Value = 853.64
$ws.Range("J99").Value = 2190.4546
$ws.Range("K99").Value = 853.64
$ws.Range("L99").Value = 2190.4546
$ws.Range("M99").Value = 644.36
$ws.Range("N99").Value = -5186.4546
$ws.Range("H105").Value = 1362922.6
$ws.Range("I105").Value = 1842412.9
$ws.Range("K105").Value = 1842412.9
$ws.Range("M105").Value = -1840665.9
$ws.Range("H134").Value = 146119.72
$ws.Range("I134").Value = 3804
$ws.Range("K134").Value = 11412
$ws.Range("M134").Value = -8877

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3128
$ws.Range("I86").Value = 3126.75
$ws.Range("K86").Value = 3126.75
$ws.Range("M86").Value = -2003.75
$ws.Range("H89").Value = 3128
$ws.Range("I89").Value = 3126.75
$ws.Range("K89").Value = 15633.75
$ws.Range("M89").Value = -10017.75
$ws.Range("H99").Value = 1334.9333
$ws.Range("I99").Value = 1337.4286
$ws.Range("K99").Value = 1337.4286
$ws.Range("M99").Value = 160.5714
$ws.Range("H105").Value = 677.86664
$ws.Range("I105").Value = 677.86664
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 677.86664
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 1069.13336
$ws.Range("N105").ClearContents()
$ws.Range("H107").Value = 482.75
$ws.Range("I107").Value = 314
$ws.Range("J107").Value = 892.5714
$ws.Range("K107").Value = 314
$ws.Range("L107").Value = 892.5714
$ws.Range("M107").Value = 1606
$ws.Range("N107").Value = -4732.5714
$ws.Range("H126").Value = 1334.9333
$ws.Range("I126").Value = 1337.4286
$ws.Range("K126").Value = 4012.2858
$ws.Range("M126").Value = -1542.2858

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 373.83334
$ws.Range("J23").Value = 414.875
$ws.Range("L23").Value = 1244.625
$ws.Range("N23").Value = -1714.625
$ws.Range("H46").Value = 2937.5757
$ws.Range("J46").Value = 2937.5757
$ws.Range("L46").Value = 8812.7271
$ws.Range("N46").Value = -8994.7271
$ws.Range("H68").Value = 2323.9314
$ws.Range("I68").Value = 862.2778
$ws.Range("J68").Value = 3746.081
$ws.Range("K68").Value = 2586.8334
$ws.Range("L68").Value = 11238.243
$ws.Range("M68").Value = -1775.8334
$ws.Range("N68").Value = -12860.243
$ws.Range("H71").Value = 2323.9314
$ws.Range("I71").Value = 862.2778
$ws.Range("J71").Value = 3746.081
$ws.Range("K71").Value = 7760.500199999999
$ws.Range("L71").Value = 33714.729
$ws.Range("M71").Value = -3704.500199999999
$ws.Range("N71").Value = -41826.729
$ws.Range("H107").Value = 702.98114
$ws.Range("J107").Value = 2100
$ws.Range("L107").Value = 6300
$ws.Range("N107").Value = -10140
$ws.Range("H109").Value = 3312.3076
$ws.Range("I109").Value = 3030
$ws.Range("J109").Value = 3335.8333
$ws.Range("K109").Value = 9090
$ws.Range("L109").Value = 10007.4999
$ws.Range("M109").Value = -8050
$ws.Range("N109").Value = -12087.4999
$ws.Range("H113").Value = 457.69696
$ws.Range("I113").Value = 512
$ws.Range("J113").Value = 384
$ws.Range("K113").Value = 1536
$ws.Range("L113").Value = 1152
$ws.Range("M113").Value = 634
$ws.Range("N113").Value = -5492
$ws.Range("H129").Value = 1537.5
$ws.Range("I129").Value = 1035.7142
$ws.Range("J129").Value = 1807.6923
$ws.Range("K129").Value = 3107.1426
$ws.Range("L129").Value = 5423.0769
$ws.Range("M129").Value = 1892.8574
$ws.Range("N129").Value = -15423.0769
$ws.Range("H131").Value = 1231.8448
$ws.Range("I131").Value = 939
$ws.Range("J131").Value = 1292.8541
$ws.Range("K131").Value = 2817
$ws.Range("L131").Value = 3878.5623
$ws.Range("M131").Value = 2223
$ws.Range("N131").Value = -13958.5623

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7671.857
$ws.Range("I80").Value = 26250
$ws.Range("J80").Value = 4575.5
$ws.Range("K80").Value = 26250
$ws.Range("L80").Value = 4575.5
$ws.Range("M80").Value = -25252
$ws.Range("N80").Value = -6571.5
$ws.Range("H83").Value = 7671.857
$ws.Range("I83").Value = 26250
$ws.Range("J83").Value = 4575.5
$ws.Range("K83").Value = 131250
$ws.Range("L83").Value = 22877.5
$ws.Range("M83").Value = -126258
$ws.Range("N83").Value = -32861.5
$ws.Range("H132").Value = 7940.8184
$ws.Range("I132").Value = 3139.3333
$ws.Range("J132").Value = 13702.6
$ws.Range("K132").Value = 9417.999899999999
$ws.Range("L132").Value = 41107.8
$ws.Range("M132").Value = -6887.999899999999
$ws.Range("N132").Value = -46167.8

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 950
$ws.Range("I22").Value = 750
$ws.Range("J22").Value = 1150
$ws.Range("K22").Value = 750
$ws.Range("L22").Value = 1150
$ws.Range("M22").Value = -455
$ws.Range("N22").Value = -1740
$ws.Range("H27").Value = 950
$ws.Range("I27").Value = 750
$ws.Range("J27").Value = 1150
$ws.Range("K27").Value = 750
$ws.Range("L27").Value = 1150
$ws.Range("M27").Value = -643
$ws.Range("N27").Value = -1364

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 729.2646999999999
$ws.Range("I113").Value = 498.2
$ws.Range("J113").Value = 769.10345
$ws.Range("K113").Value = 1494.6
$ws.Range("L113").Value = 2307.31035
$ws.Range("M113").Value = 675.4000000000001
$ws.Range("N113").Value = -6647.31035
$ws.Range("H132").Value = 3984.6
$ws.Range("I132").Value = 3905.3845
$ws.Range("J132").Value = 4499.5
$ws.Range("K132").Value = 11716.1535
$ws.Range("L132").Value = 13498.5
$ws.Range("M132").Value = -9186.1535
$ws.Range("N132").Value = -18558.5
$ws.Range("H135").Value = 57500
$ws.Range("J135").Value = 57500
$ws.Range("L135").Value = 57500
$ws.Range("N135").Value = -67640
